# who_test_controls_analysis_results.xlsx edit script
# - Clears stray literal "None" text left over from the old string-formatting
#   bug in the Enhancement Feedback sheet (debugging the token/format error).
# - Fixes a corrupted entry on the Multi-Control Candidates sheet (re-adds the
#   leading apostrophe that was being swallowed) and rebalances its column
#   widths.
# - Removes accidental bold styling from two headings and corrects the
#   "seven key elements" -> "five key elements" copy on the Methodology sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Enhancement Feedback" sheet: blank out cells that literally contain the
#    word "None" (the WHO/WHEN/WHAT/ESCALATION feedback columns), leaving the
#    WHY column (E) untouched.
# ---------------------------------------------------------------------------
$wsFeedback = $wb.Worksheets.Item("Enhancement Feedback")
$feedbackCols = @("B", "C", "D", "F")
for ($row = 2; $row -le 9; $row++) {
    foreach ($col in $feedbackCols) {
        $cell = $wsFeedback.Range("$col$row")
        if ($cell.Text -eq "None") {
            $cell.Value = ""
        }
    }
}

# ---------------------------------------------------------------------------
# 2) "Multi-Control Candidates" sheet: repair the truncated/corrupted text in
#    B3 (it lost its leading apostrophe) and resize the columns.
# ---------------------------------------------------------------------------
$wsMulti = $wb.Worksheets.Item("Multi-Control Candidates")

$b3 = $wsMulti.Range("B3")
$b3.Value = "''-facing materials must be reviewed"
$b3.ClearFormats()

$wsMulti.Columns.Item(1).ColumnWidth = 12.366666666666665
$wsMulti.Columns.Item(2).ColumnWidth = 49.166666666666664
$wsMulti.Columns.Item(3).ColumnWidth = 6.866666666666668
$wsMulti.Columns.Item(4).ColumnWidth = 49.166666666666664

# ---------------------------------------------------------------------------
# 3) "Executive Summary" sheet: the "Top Vague Terms" heading should not be
#    bold anymore.
# ---------------------------------------------------------------------------
$wsExec = $wb.Worksheets.Item("Executive Summary")
$wsExec.Range("A20").ClearFormats()

# ---------------------------------------------------------------------------
# 4) "Methodology" sheet: correct the element count in the overview copy and
#    remove the stray bold styling from the "Enhanced Validation Checks"
#    heading.
# ---------------------------------------------------------------------------
$wsMethod = $wb.Worksheets.Item("Methodology")
$wsMethod.Range("A4").Value = "This analysis evaluates control descriptions based on five key elements that should be present in a well-written control description:"
$wsMethod.Range("A42").ClearFormats()
